$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("teaching")

# Clear the "with" column (C) values for the data rows - no longer used
$ws.Range("C2").Value = $null
$ws.Range("C3").Value = $null
$ws.Range("C4").Value = $null

# Update "what" column (D) content for each row
$ws.Range("D2").Value = "Behavioural ecology and sociobiology; BMC biology; Insects; Journal of Evolutionary Biology, Phil. Transactions of the Royal Society; PLoS Genetics; Proceedings of the Royal Society; Molecular Ecology; Frontiers in BioEngineering"
$ws.Range("D3").Value = "MRC fellowships; BBSRC Fellowships; GWIS National Fellowships"
$ws.Range("D4").Value = "Fellow of the Genetics Societ; Fellow of the Royal Statistical Society; Vectorbite; National Network of Bioscience Educators; Advance HE"

# Update the selection to match the saved cursor position
$ws.Range("D4").Select()
